$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '68.675.78'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '3.852.39'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '524.79'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +7.55%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.53'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.606'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.40%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.711'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -4.31%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.169'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -5.50%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000326'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -7.16%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '41.66'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.05%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '4.458.19'
$ws.Range('E13').Value = '  -1.43%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.16'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.07%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.51'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +7.55%  '
$ws.Range('D16').Value = '3.832.49'
$ws.Range('E16').Value = '  -1.91%  '
$ws.Range('E17').Value = '  +7.12%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '13.90'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.11%  '
$ws.Range('E19').Value = '  -1.53%  '
$ws.Range('D20').Value = '68.665.49'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '418.93'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.36%  '
$ws.Range('E22').Value = '  -3.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.02'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -4.61%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '86.87'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.15%  '
$ws.Range('E25').Value = '  +5.72%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.30'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -9.70%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.56'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -4.26%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '35.94'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -3.55%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '680.76'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.44%  '
$ws.Range('E30').Value = '  -2.29%  '
$ws.Range('E31').Value = '  -3.64%  '
$ws.Range('E32').Value = '  -2.93%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '66.98'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +9.26%  '
$ws.Range('E34').Value = '  +1.99%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.87'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -4.72%  '
$ws.Range('B36').Value = 'PEPE'
$ws.Range('C36').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D36').Value = '0.0₃0847'
$ws.Range('E36').Value = '  -3.44%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '39.78'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.78%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('E39').Value = '  -1.57%  '
$ws.Range('E40').Value = '  -0.22%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0477'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.22%  '
$ws.Range('E42').Value = '  +1.51%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.15'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.40%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.72'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -10.87%  '
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('E46').Value = '  -2.59%  '
$ws.Range('D47').Value = '2.768.76'
$ws.Range('E47').Value = '  +14.60%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.94'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +4.80%  '
$ws.Range('E49').Value = '  +12.00%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '144.25'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.30%  '
$ws.Range('E51').Value = '  -3.26%  '
